# Apply the latest crypto price/volume snapshot pulled by the scraper.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '65.710.54'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.78%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.445.10'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -2.77%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.05%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '592.67'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.63%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '138.27'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -5.14%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.442.96'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.89%  '

# Row 8
$ws.Range("E8").Value = '  +0.10%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.507'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.37%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.37'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -5.15%  '

# Row 11
$ws.Range("E11").Value = '  -7.64%  '

# Row 12
$ws.Range("E12").Value = '  -6.50%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.030.40'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.69%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000182'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -9.14%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '26.63'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -8.17%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.444.64'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.97%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '65.616.82'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.82%  '

# Row 18
$ws.Range("E18").Value = '  -1.60%  '

# Row 19
$ws.Range("E19").Value = '  -10.48%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.93'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -5.05%  '

# Row 21
$ws.Range("E21").Value = '  -5.68%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '396.03'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -5.11%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.556'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -7.32%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '73.70'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -5.05%  '

# Row 25
$ws.Range("E25").Value = '  -0.02%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.589.19'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.49%  '

# Row 27
$ws.Range("E27").Value = '  -7.12%  '

# Row 28
$ws.Range("E28").Value = '  +0.00%  '

# Row 29
$ws.Range("E29").Value = '  -8.23%  '

# Row 30
$ws.Range("E30").Value = '  -7.34%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.25'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -8.70%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.453.86'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.41%  '

# Row 33
$ws.Range("E33").Value = '  +0.01%  '

# Row 34
$ws.Range("E34").Value = '  -6.96%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '23.09'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -5.37%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '172.96'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.48%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.97'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -7.76%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.20'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -7.63%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.49'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -7.37%  '

# Row 40
$ws.Range("E40").Value = '  -8.11%  '

# Row 41
$ws.Range("E41").Value = '  -6.07%  '

# Row 42
$ws.Range("E42").Value = '  -3.84%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '43.82'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.98%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.00'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.07%  '

# Row 45
$ws.Range("E45").Value = '  -12.35%  '

# Row 46
$ws.Range("E46").Value = '  -9.43%  '

# Row 47
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '23.21'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.00%  '

# Row 48
$ws.Range("B48").Value = 'ONDO'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.12'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.01%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.61'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -6.94%  '

# Row 50
$ws.Range("E50").Value = '  -11.89%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.213.72'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -6.62%  '
